# Apply updates to "ASR Results" sheet per commit:
# "F04 Froze Token Embeddings and Decoder 12"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B text updates (inline strings wrapped in <...>)
$ws.Range("B5").Value  = "<foot>"
$ws.Range("B10").Value = "<eis>"
$ws.Range("B12").Value = "<their>"
$ws.Range("B13").Value = "<sera>"
$ws.Range("B14").Value = "<al>"
$ws.Range("B17").Value = "<seen>"
$ws.Range("B18").Value = "<tab>"

# Column C numeric updates
$ws.Range("C2").Value  = 33
$ws.Range("C3").Value  = 33
$ws.Range("C4").Value  = 24
$ws.Range("C5").Value  = 35
$ws.Range("C6").Value  = 25
$ws.Range("C7").Value  = 26
$ws.Range("C8").Value  = 27
$ws.Range("C9").Value  = 34
$ws.Range("C10").Value = 31
$ws.Range("C11").Value = 32
$ws.Range("C12").Value = 43
$ws.Range("C14").Value = 22
$ws.Range("C15").Value = 31
$ws.Range("C16").Value = 27
$ws.Range("C17").Value = 29
$ws.Range("C18").Value = 23
